# Apply updated crypto market data (price + 1h volume change) scraped by the
# GitHub Actions job. Also reflects the new ranking order for rows 33-34
# (Fetch.AI now ranks above Aptos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cell updates: cell reference -> new text value.
$textUpdates = @{
    "D2" = "56.354.30"
    "E2" = "  -2.20%  "
    "D3" = "2.376.27"
    "E3" = "  -1.62%  "
    "E4" = "  -0.10%  "
    "E5" = "  -1.35%  "
    "E6" = "  -2.08%  "
    "E7" = "  +0.50%  "
    "E8" = "  -2.57%  "
    "D9" = "2.381.21"
    "E9" = "  -2.94%  "
    "E10" = "  -0.26%  "
    "E11" = "  +0.73%  "
    "E12" = "  +0.69%  "
    "E13" = "  +0.36%  "
    "D14" = "2.794.95"
    "E14" = "  -1.90%  "
    "D15" = "56.289.40"
    "E15" = "  -1.97%  "
    "E16" = "  -1.58%  "
    "E17" = "  -1.74%  "
    "D18" = "2.385.42"
    "E18" = "  -2.01%  "
    "E19" = "  -2.79%  "
    "E20" = "  -2.92%  "
    "E21" = "  -2.57%  "
    "E22" = "  -2.10%  "
    "E24" = "  -1.26%  "
    "E25" = "  +0.33%  "
    "E26" = "  -3.91%  "
    "E27" = "  -4.67%  "
    "E28" = "  -4.22%  "
    "E29" = "  -0.92%  "
    "D30" = "0.0₃0715"
    "E30" = "  -3.42%  "
    "E31" = "  -3.54%  "
    "E32" = "  +0.25%  "
    "B33" = "Fetch.AI"
    "C33" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "E33" = "  -4.98%  "
    "B34" = "Aptos"
    "C34" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "E34" = "  -7.71%  "
    "E35" = "  +0.51%  "
    "E36" = "  -2.69%  "
    "E37" = "  -6.18%  "
    "E38" = "  -3.11%  "
    "E39" = "  -1.41%  "
    "E40" = "  -4.19%  "
    "E41" = "  -5.13%  "
    "E42" = "  -2.77%  "
    "E43" = "  -2.12%  "
    "E44" = "  -5.50%  "
    "E45" = "  -1.17%  "
    "E46" = "  -1.74%  "
    "E47" = "  -7.43%  "
    "E48" = "  -2.73%  "
    "E49" = "  -2.94%  "
    "E50" = "  -1.17%  "
    "E51" = "  -3.60%  "
}

# Price cells whose new text would otherwise be auto-parsed by Excel as a
# number (and so lose trailing zeros / become a float) - force text format.
$numericLookingUpdates = @{
    "D5" = "502.07"
    "D6" = "130.24"
    "D7" = "0.998"
    "D8" = "0.545"
    "D10" = "0.0984"
    "D12" = "0.325"
    "D13" = "4.66"
    "D16" = "21.63"
    "D19" = "10.07"
    "D20" = "4.02"
    "D21" = "307.41"
    "D24" = "64.64"
    "D28" = "7.33"
    "D29" = "172.51"
    "D31" = "1.64"
    "D33" = "1.09"
    "D34" = "5.76"
    "D35" = "0.997"
    "D36" = "17.56"
    "D38" = "3.78"
    "D39" = "36.01"
    "D40" = "0.793"
    "D41" = "1.40"
    "D42" = "130.98"
    "D43" = "3.36"
    "D44" = "4.75"
    "D46" = "0.0904"
    "D47" = "241.31"
    "D48" = "0.0484"
    "D50" = "17.06"
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
    $cell.Style = "Normal"
}
